# Add two new worksheets ("selectFlight" and "bookFlight") right after
# "openBrowser", populate them with flight-selection data, and move the
# selection on "inputFlightDetails" to B6.

$wb = $excel.ActiveWorkbook

# --- Set the selection on the existing "inputFlightDetails" sheet first so
# that the later sheet-activation (from adding new sheets) ends up being the
# one that Excel reports as the active tab.
$flightSheet = $wb.Worksheets.Item("inputFlightDetails")
$flightSheet.Range("B6").Select()

# --- Insert "selectFlight" right after "openBrowser" ------------------------
$openBrowser = $wb.Worksheets.Item("openBrowser")
$selectFlight = $wb.Worksheets.Add($null, $openBrowser)
$selectFlight.Name = "selectFlight"

$selectFlight.Range("A3").Value = "Departure Flight"
$selectFlight.Range("B3").Value = "Pangaea Airlines 362"
$selectFlight.Range("A4").Value = "Return Flight"
$selectFlight.Range("B4").Value = "Unified Airlines 633"
$selectFlight.Range("A1").Value = "Departure"
$selectFlight.Range("B1").Value = "Acapulco"
$selectFlight.Range("A2").Value = "Arrival"
$selectFlight.Range("B2").Value = "Zurich"
$selectFlight.Range("A1:B4").Select()

# --- Insert "bookFlight" right after "selectFlight" -------------------------
$bookFlight = $wb.Worksheets.Add($null, $selectFlight)
$bookFlight.Name = "bookFlight"

$bookFlight.Range("A3").Value = "Departure Flight"
$bookFlight.Range("B3").Value = "Pangaea Airlines 362"
$bookFlight.Range("A4").Value = "Return Flight"
$bookFlight.Range("B4").Value = "Unified Airlines 633"
$bookFlight.Range("A1").Value = "Departure"
$bookFlight.Range("B1").Value = "Acapulco"
$bookFlight.Range("A2").Value = "Arrival"
$bookFlight.Range("B2").Value = "Zurich"
$bookFlight.Range("A1:B4").Select()

$bookFlight.Activate()
